$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eCL Jobs")

# ---------------------------------------------------------------------------
# Helper: move the 5-cell block A:E of row $src down to row $dst, carrying
# both values and cell formatting, plus the row height (only meaningful for
# data rows - separator rows keep default height).
# ---------------------------------------------------------------------------
function Move-JobRow($src, $dst) {
    $srcRange = $ws.Range("A" + $src + ":E" + $src)
    $dstRange = $ws.Range("A" + $dst + ":E" + $dst)

    $srcRange.Copy()
    $dstRange.PasteSpecial(-4163)   # xlPasteValues
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)   # xlPasteFormats

    $srcHeight = $ws.Rows($src).RowHeight
    if ($srcHeight -ne 15) {
        $ws.Rows($dst).RowHeight = $srcHeight
    }
}

# ---------------------------------------------------------------------------
# Phase 1: shift the existing 10 job blocks (data row + blank separator row)
# down to their new target rows so room is made for the 4 new job blocks.
# Processed from the bottom (largest destination) upward so that a source
# row is never clobbered before it has been read.
# ---------------------------------------------------------------------------
Move-JobRow 20 28
Move-JobRow 19 23
Move-JobRow 18 22
Move-JobRow 17 21
Move-JobRow 16 20
Move-JobRow 15 19
Move-JobRow 14 18
Move-JobRow 13 15
Move-JobRow 12 14
Move-JobRow 11 13
Move-JobRow 10 12
Move-JobRow 9  11
Move-JobRow 8  10
Move-JobRow 7  9
Move-JobRow 6  8

# The separator row that used to sit at 19 (now living at row 23) had the
# "heavy" wrap-text separator style; in the new layout it reverts back to
# the plain separator style (same as row 9/13/15/17/19/21).
$ws.Range("C13:E13").Copy()
$ws.Range("C23:E23").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# Phase 2: populate the newly-freed rows with the four new job blocks, each
# formatted like the standard data-row / separator-row pair used throughout
# the sheet (copied from an existing, untouched pair: rows 4 and 5).
# ---------------------------------------------------------------------------
function New-JobRow($row, $height, $a, $b, $c, $d, $e) {
    $tmplData = $ws.Range("A4:E4")
    $dstData = $ws.Range("A" + $row + ":E" + $row)
    $tmplData.Copy()
    $dstData.PasteSpecial(-4122)   # xlPasteFormats
    $ws.Rows($row).RowHeight = $height

    $ws.Range("A" + $row).Value2 = $a
    $ws.Range("B" + $row).Value2 = $b
    $ws.Range("C" + $row).Value2 = $c
    $ws.Range("D" + $row).Value2 = $d
    $ws.Range("E" + $row).Value2 = $e

    $sepRow = $row + 1
    $tmplSep = $ws.Range("A5:E5")
    $dstSep = $ws.Range("A" + $sepRow + ":E" + $sepRow)
    $tmplSep.Copy()
    $dstSep.PasteSpecial(-4122)    # xlPasteFormats
    $dstSep.ClearContents()
}

New-JobRow 6 75 "CoachingGenericLoad" "Generic file Load" "Inports generic coaching logs" "\\vrivscors01\BCC Scorecards\Coaching\Generic\eCL_Generic_Feed_XXX<YYYYMMDD>.csv" "EC.Generic_Coaching_Stage`n    EC.Generic_Coaching_Rejected`n    EC.Generic_Coaching_Fact`n        EC.Coaching_Log`n        EC.Coaching_Log_Reason"

New-JobRow 16 75 "CoachingQualityOtherLoad" "Quality Other Load" "Imports quality reports information to create coaching logs" "\\vrivscors01\BCC Scorecards\Coaching\Quality\eCL_Quality_Feed_XXX<YYYYMMDD>.csv" "EC.Quality_Other_Coaching_Stage`n    EC.Quality_Other_Coacing_Rejected`n    EC.Quality_Other_Coaching_Fact`n        EC.Coaching_Log`n        EC.Coaching_Log_Reason"

New-JobRow 24 90 "CoachingSurveyReminders" "SurveyReminders" "Sends reminder email notification regarding survey" "EC.Survey_Response_Header.NotificationDate = today's date - 3" "EC.Survey_Response_Header.ReminderSent = 1`nEC.Survey_Response_Header.ReminderDate = email date`nEC.Survey_Response_Header.ReminderCount = +1`nemail reminder notification sent to recipient"

New-JobRow 26 75 "CoachingTraining" "Load Training Files" "Imports training reports information to create coaching logs" "\\vrivscors01\BCC Scorecards\Coaching\Training\eCL_Training_Feed_XXX<YYYYMMDD>.csv" "EC.Training_Coaching_Stage`n    EC.Training_Coacing_Rejected`n    EC.Training_Coaching_Fact`n        EC.Coaching_Log`n        EC.Coaching_Log_Reason"

# Rows 27 and 29 (the separators that trail the newly-added "CoachingTraining"
# and "CoachingWHLoad" blocks respectively) use the "heavy" wrap-text
# separator style, matching row 11 (untouched throughout this script).
$tmplHeavySep = $ws.Range("A11:E11")
foreach ($r in 27, 29) {
    $dst = $ws.Range("A" + $r + ":E" + $r)
    $tmplHeavySep.Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats
    $dst.ClearContents()
}

$ws.Range("A1").Select()
